$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''27.887.22'

# Row 3
$ws.Range('D3').Value = '''1.880.46'
$ws.Range('E3').Value = '''  +1.52%  '

# Row 4
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '''  +0.13%  '

# Row 5
$ws.Range('D5').Value = '''332.91'
$ws.Range('E5').Value = '''  +3.18%  '

# Row 6
$ws.Range('D6').Value = '''1.002'
$ws.Range('E6').Value = '''  +0.11%  '

# Row 7
$ws.Range('D7').Value = '''0.4739'
$ws.Range('E7').Value = '''  +5.40%  '

# Row 8
$ws.Range('D8').Value = '''0.3972'
$ws.Range('E8').Value = '''  +3.67%  '

# Row 9
$ws.Range('D9').Value = '''48.13'
$ws.Range('E9').Value = '''  -0.35%  '

# Row 10
$ws.Range('D10').Value = '''0.08058'
$ws.Range('E10').Value = '''  +2.61%  '

# Row 11
$ws.Range('E11').Value = '''  +1.44%  '

# Row 12
$ws.Range('E12').Value = '''  +2.65%  '

# Row 13
$ws.Range('D13').Value = '''1.897.03'
$ws.Range('E13').Value = '''  +3.62%  '

# Row 14
$ws.Range('D14').Value = '''5.966'
$ws.Range('E14').Value = '''  +1.91%  '

# Row 15
$ws.Range('D15').Value = '''7.213'
$ws.Range('E15').Value = '''  +1.20%  '

# Row 16
$ws.Range('D16').Value = '''1.002'
$ws.Range('E16').Value = '''  +0.11%  '

# Row 17
$ws.Range('D17').Value = '''0.00001052'
$ws.Range('E17').Value = '''  +2.18%  '

# Row 18
$ws.Range('D18').Value = '''87.27'
$ws.Range('E18').Value = '''  +1.81%  '

# Row 19
$ws.Range('D19').Value = '''0.06624'
$ws.Range('E19').Value = '''  +1.95%  '

# Row 20
$ws.Range('E20').Value = '''  +1.81%  '

# Row 21
$ws.Range('E21').Value = '''  +0.14%  '

# Row 22
$ws.Range('D22').Value = '''27.978.41'
$ws.Range('E22').Value = '''  +2.43%  '

# Row 23
$ws.Range('D23').Value = '''5.516'
$ws.Range('E23').Value = '''  +0.86%  '

# Row 24
$ws.Range('E24').Value = '''  +2.67%  '

# Row 25
$ws.Range('D25').Value = '''2.310'
$ws.Range('E25').Value = '''  +2.11%  '

# Row 26
$ws.Range('D26').Value = '''2.123.36'
$ws.Range('E26').Value = '''  +3.38%  '

# Row 27
$ws.Range('D27').Value = '''157.94'
$ws.Range('E27').Value = '''  +4.25%  '

# Row 28
$ws.Range('D28').Value = '''20.28'
$ws.Range('E28').Value = '''  +4.63%  '

# Row 29
$ws.Range('E29').Value = '''  +2.75%  '

# Row 30
$ws.Range('D30').Value = '''5.629'
$ws.Range('E30').Value = '''  +1.56%  '

# Row 31
$ws.Range('D31').Value = '''122.69'
$ws.Range('E31').Value = '''  +2.38%  '

# Row 32
$ws.Range('D32').Value = '''0.9874'
$ws.Range('E32').Value = '''  +5.73%  '

# Row 33
$ws.Range('D33').Value = '''0.09585'
$ws.Range('E33').Value = '''  +2.85%  '

# Row 34
$ws.Range('D34').Value = '''1.468'
$ws.Range('E34').Value = '''  -0.48%  '

# Row 35
$ws.Range('D35').Value = '''3.620'
$ws.Range('E35').Value = '''  +0.73%  '

# Row 36
$ws.Range('D36').Value = '''5.337'
$ws.Range('E36').Value = '''  +1.50%  '

# Row 37
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.02273'
$ws.Range('E37').Value = '''  +2.28%  '

# Row 38
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '''0.06126'
$ws.Range('E38').Value = '''  +2.53%  '

# Row 39
$ws.Range('D39').Value = '''1.241'
$ws.Range('E39').Value = '''  +3.27%  '

# Row 40
$ws.Range('D40').Value = '''8.257'
$ws.Range('E40').Value = '''  -0.23%  '

# Row 41
$ws.Range('D41').Value = '''0.6037'
$ws.Range('E41').Value = '''  +2.45%  '

# Row 42
$ws.Range('E42').Value = '''  +0.16%  '

# Row 43
$ws.Range('D43').Value = '''0.1904'
$ws.Range('E43').Value = '''  +2.88%  '

# Row 44
$ws.Range('E44').Value = '''  +1.79%  '

# Row 45
$ws.Range('D45').Value = '''1.264'
$ws.Range('E45').Value = '''  +1.01%  '

# Row 46
$ws.Range('D46').Value = '''0.5721'
$ws.Range('E46').Value = '''  +1.44%  '

# Row 47
$ws.Range('D47').Value = '''12.34'
$ws.Range('E47').Value = '''  +1.32%  '

# Row 48
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''1.952'
$ws.Range('E48').Value = '''  +1.34%  '

# Row 49
$ws.Range('B49').Value = 'PancakeSwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D49').Value = '''3.412'
$ws.Range('E49').Value = '''  +1.53%  '

# Row 50
$ws.Range('D50').Value = '''0.06838'
$ws.Range('E50').Value = '''  -0.53%  '

# Row 51
$ws.Range('D51').Value = '''113.78'
$ws.Range('E51').Value = '''  +5.10%  '
